# Applies the cryptos-list price/volume refresh described in the commit.
# Numeric-looking Price (column D) values are written through a text-format
# round-trip so Excel keeps them as literal strings (matching the source data,
# which stores every Price/Volume cell as text) instead of auto-coercing them
# to numbers; the temporary style is restored immediately afterwards so no
# cell ends up with a different style than it started with.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Value = "63.639.86"; Numeric = $false },
    @{ Addr = "E2"; Value = "  +2.78%  "; Numeric = $false },
    @{ Addr = "D3"; Value = "3.127.54"; Numeric = $false },
    @{ Addr = "E3"; Value = "  +1.60%  "; Numeric = $false },
    @{ Addr = "E4"; Value = "  +0.12%  "; Numeric = $false },
    @{ Addr = "D5"; Value = "587.06"; Numeric = $true },
    @{ Addr = "E5"; Value = "  +1.21%  "; Numeric = $false },
    @{ Addr = "D6"; Value = "146.21"; Numeric = $true },
    @{ Addr = "E6"; Value = "  +2.81%  "; Numeric = $false },
    @{ Addr = "E7"; Value = "  +0.08%  "; Numeric = $false },
    @{ Addr = "D8"; Value = "3.117.45"; Numeric = $false },
    @{ Addr = "E8"; Value = "  +1.58%  "; Numeric = $false },
    @{ Addr = "D9"; Value = "0.532"; Numeric = $true },
    @{ Addr = "E9"; Value = "  +0.98%  "; Numeric = $false },
    @{ Addr = "D10"; Value = "0.159"; Numeric = $true },
    @{ Addr = "E10"; Value = "  +13.29%  "; Numeric = $false },
    @{ Addr = "D11"; Value = "5.69"; Numeric = $true },
    @{ Addr = "E11"; Value = "  +2.18%  "; Numeric = $false },
    @{ Addr = "D12"; Value = "0.468"; Numeric = $true },
    @{ Addr = "E12"; Value = "  +0.05%  "; Numeric = $false },
    @{ Addr = "D13"; Value = "0.0000250"; Numeric = $true },
    @{ Addr = "E13"; Value = "  +4.28%  "; Numeric = $false },
    @{ Addr = "D14"; Value = "36.37"; Numeric = $true },
    @{ Addr = "E14"; Value = "  +2.86%  "; Numeric = $false },
    @{ Addr = "E15"; Value = "  -0.62%  "; Numeric = $false },
    @{ Addr = "D16"; Value = "3.646.41"; Numeric = $false },
    @{ Addr = "E16"; Value = "  +1.64%  "; Numeric = $false },
    @{ Addr = "D17"; Value = "7.16"; Numeric = $true },
    @{ Addr = "E17"; Value = "  -1.75%  "; Numeric = $false },
    @{ Addr = "D18"; Value = "63.573.57"; Numeric = $false },
    @{ Addr = "E18"; Value = "  +2.82%  "; Numeric = $false },
    @{ Addr = "D19"; Value = "3.127.29"; Numeric = $false },
    @{ Addr = "E19"; Value = "  +1.69%  "; Numeric = $false },
    @{ Addr = "D20"; Value = "461.74"; Numeric = $true },
    @{ Addr = "E20"; Value = "  +2.69%  "; Numeric = $false },
    @{ Addr = "D21"; Value = "14.38"; Numeric = $true },
    @{ Addr = "E21"; Value = "  +3.29%  "; Numeric = $false },
    @{ Addr = "D22"; Value = "0.732"; Numeric = $true },
    @{ Addr = "E22"; Value = "  +0.40%  "; Numeric = $false },
    @{ Addr = "D23"; Value = "7.52"; Numeric = $true },
    @{ Addr = "E23"; Value = "  +1.26%  "; Numeric = $false },
    @{ Addr = "D24"; Value = "13.21"; Numeric = $true },
    @{ Addr = "E24"; Value = "  -4.05%  "; Numeric = $false },
    @{ Addr = "D25"; Value = "82.13"; Numeric = $true },
    @{ Addr = "E25"; Value = "  +0.27%  "; Numeric = $false },
    @{ Addr = "D27"; Value = "8.89"; Numeric = $true },
    @{ Addr = "E27"; Value = "  +9.02%  "; Numeric = $false },
    @{ Addr = "E28"; Value = "  +1.13%  "; Numeric = $false },
    @{ Addr = "D29"; Value = "2.22"; Numeric = $true },
    @{ Addr = "E29"; Value = "  -2.18%  "; Numeric = $false },
    @{ Addr = "D31"; Value = "6.86"; Numeric = $true },
    @{ Addr = "E31"; Value = "  +1.41%  "; Numeric = $false },
    @{ Addr = "D32"; Value = "27.02"; Numeric = $true },
    @{ Addr = "E32"; Value = "  +0.92%  "; Numeric = $false },
    @{ Addr = "E33"; Value = "  -1.78%  "; Numeric = $false },
    @{ Addr = "D34"; Value = "0.0₃0862"; Numeric = $false },
    @{ Addr = "E34"; Value = "  +7.94%  "; Numeric = $false },
    @{ Addr = "D35"; Value = "2.36"; Numeric = $true },
    @{ Addr = "E35"; Value = "  +7.31%  "; Numeric = $false },
    @{ Addr = "E36"; Value = "  +1.08%  "; Numeric = $false },
    @{ Addr = "D37"; Value = "3.33"; Numeric = $true },
    @{ Addr = "E37"; Value = "  +11.17%  "; Numeric = $false },
    @{ Addr = "D38"; Value = "6.07"; Numeric = $true },
    @{ Addr = "E38"; Value = "  +0.22%  "; Numeric = $false },
    @{ Addr = "D39"; Value = "50.92"; Numeric = $true },
    @{ Addr = "E39"; Value = "  +1.06%  "; Numeric = $false },
    @{ Addr = "D40"; Value = "447.00"; Numeric = $true },
    @{ Addr = "E40"; Value = "  +3.95%  "; Numeric = $false },
    @{ Addr = "D41"; Value = "8.73"; Numeric = $true },
    @{ Addr = "E41"; Value = "  -1.11%  "; Numeric = $false },
    @{ Addr = "D42"; Value = "0.0371"; Numeric = $true },
    @{ Addr = "E42"; Value = "  -0.10%  "; Numeric = $false },
    @{ Addr = "D43"; Value = "2.889.41"; Numeric = $false },
    @{ Addr = "E43"; Value = "  +3.26%  "; Numeric = $false },
    @{ Addr = "D44"; Value = "0.279"; Numeric = $true },
    @{ Addr = "E44"; Value = "  +1.99%  "; Numeric = $false },
    @{ Addr = "E45"; Value = "  +1.64%  "; Numeric = $false },
    @{ Addr = "D46"; Value = "2.17"; Numeric = $true },
    @{ Addr = "E46"; Value = "  +2.41%  "; Numeric = $false },
    @{ Addr = "D47"; Value = "36.16"; Numeric = $true },
    @{ Addr = "E47"; Value = "  +3.02%  "; Numeric = $false },
    @{ Addr = "B48"; Value = "USDe"; Numeric = $false },
    @{ Addr = "C48"; Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; Numeric = $false },
    @{ Addr = "D48"; Value = "0.999"; Numeric = $true },
    @{ Addr = "E48"; Value = "  +0.06%  "; Numeric = $false },
    @{ Addr = "B49"; Value = "Monero"; Numeric = $false },
    @{ Addr = "C49"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; Numeric = $false },
    @{ Addr = "D49"; Value = "124.52"; Numeric = $true },
    @{ Addr = "E49"; Value = "  +0.65%  "; Numeric = $false },
    @{ Addr = "E50"; Value = "  -0.23%  "; Numeric = $false },
    @{ Addr = "D51"; Value = "24.58"; Numeric = $true },
    @{ Addr = "E51"; Value = "  +2.05%  "; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.Numeric) {
        # Force text storage: switch to a Text format, assign the literal
        # string, then restore the cell's original style so formatting is
        # unaffected by the round-trip.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $u.Value
    }
}
